# Refresh the crypto price/volume table (GitHub Actions data pull).
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). Prices/percentages are
# stored as text (Excel can't natively hold "60.219.65"-style grouped
# numbers), so for plain decimal-looking prices we force the cell to Text
# format first to stop Excel's COM layer from silently re-interpreting the
# assignment as a real number (e.g. "1.00" -> 1, "0.0374" -> 3.74E-02).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.219.65'
$ws.Range('E2').Value = '  -2.77%  '
$ws.Range('D3').Value = '2.943.13'
$ws.Range('E3').Value = '  -2.62%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '519.63'
$ws.Range('E5').Value = '  -1.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '128.52'
$ws.Range('E6').Value = '  -0.39%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '2.938.80'
$ws.Range('E8').Value = '  -2.51%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.477'
$ws.Range('E9').Value = '  -1.88%  '
$ws.Range('E10').Value = '  +2.23%  '
$ws.Range('E11').Value = '  -1.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.433'
$ws.Range('E12').Value = '  -2.41%  '
$ws.Range('E13').Value = '  -1.73%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.68'
$ws.Range('E14').Value = '  -1.61%  '
$ws.Range('D15').Value = '3.427.42'
$ws.Range('E15').Value = '  -1.49%  '
$ws.Range('E16').Value = '  -0.03%  '
$ws.Range('D17').Value = '60.312.23'
$ws.Range('E17').Value = '  -2.71%  '
$ws.Range('D18').Value = '2.948.39'
$ws.Range('E18').Value = '  -2.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.40'
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '452.49'
$ws.Range('E20').Value = '  -3.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.89'
$ws.Range('E21').Value = '  -0.56%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.662'
$ws.Range('E22').Value = '  -2.87%  '
$ws.Range('E23').Value = '  -3.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.33'
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.53'
$ws.Range('E26').Value = '  -1.24%  '
$ws.Range('E27').Value = '  -0.91%  '
$ws.Range('E28').Value = '  -5.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '25.01'
$ws.Range('E30').Value = '  -0.94%  '
$ws.Range('E31').Value = '  +2.99%  '
$ws.Range('E32').Value = '  -0.27%  '
$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '54.22'
$ws.Range('E33').Value = '  -3.00%  '
$ws.Range('B34').Value = 'Stacks'
$ws.Range('C34').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.23'
$ws.Range('E34').Value = '  -4.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.24'
$ws.Range('E35').Value = '  +2.58%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.67'
$ws.Range('E36').Value = '  -1.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '443.82'
$ws.Range('E37').Value = '  -3.06%  '
$ws.Range('D38').Value = '3.128.72'
$ws.Range('E38').Value = '  +2.93%  '
$ws.Range('E39').Value = '  -0.56%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0374'
$ws.Range('E40').Value = '  -2.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.115'
$ws.Range('E41').Value = '  +3.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.88'
$ws.Range('E42').Value = '  +0.30%  '
$ws.Range('E43').Value = '  -3.86%  '
$ws.Range('E45').Value = '  -0.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '24.72'
$ws.Range('E46').Value = '  +4.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '118.20'
$ws.Range('E47').Value = '  +2.86%  '
$ws.Range('E48').Value = '  +0.54%  '
$ws.Range('B49').Value = 'PEPE'
$ws.Range('C49').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D49').Value = '0.0₃0496'
$ws.Range('E49').Value = '  -2.64%  '
$ws.Range('B50').Value = 'Fetch.AI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.89'
$ws.Range('E50').Value = '  -3.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.23'
$ws.Range('E51').Value = '  +6.60%  '
